$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 3: "Lack of experience in weather prediction algorithms" row.
# Content is unchanged; only the row height grows (auto little bit taller).
# ---------------------------------------------------------------------------
$ws.Rows.Item(3).RowHeight = 72.75

# ---------------------------------------------------------------------------
# Row 4: "Missing a project meeting" row - risk text re-cased.
# ---------------------------------------------------------------------------
$ws.Range("C4").Value = "Missing a project meeting "

# ---------------------------------------------------------------------------
# Row 5: "Programming issues" row - risk text re-cased.
# ---------------------------------------------------------------------------
$ws.Range("C5").Value = "Programming issues "

# ---------------------------------------------------------------------------
# New row 6: "Lack of C# experience"
# ---------------------------------------------------------------------------
$ws.Range("C6").Value = "Lack of C# experience"
$ws.Range("D6").Value = "Medium "
$ws.Range("E6").Value = "Medium "
$ws.Range("F6").Value = "Increased time on tasks. Or increase workload for those with C# experience"
$ws.Range("G6").Value = "Carry out training so the whole team has a basic understanidng of the programming language "
$ws.Range("H6").Value = "Low "
$ws.Range("I6").Value = "Low "

# ---------------------------------------------------------------------------
# New row 7: "Little experience in database building/ management"
# ---------------------------------------------------------------------------
$ws.Range("C7").Value = "Little experience in database building/ management"
$ws.Range("D7").Value = "High "
$ws.Range("E7").Value = "Medium "
$ws.Range("F7").Value = "An unsuitable database will be created that will store weather and user data incorrectly"
$ws.Range("G7").Value = "Conduct research on database building"
$ws.Range("H7").Value = "Medium "
$ws.Range("I7").Value = "Low "

# ---------------------------------------------------------------------------
# New row 8: "Little experience in server development"
# ---------------------------------------------------------------------------
$ws.Range("C8").Value = "Little experience in server development"
$ws.Range("D8").Value = "High "
$ws.Range("E8").Value = "Medium "
$ws.Range("F8").Value = "We may struggle to successfully store our database  "
$ws.Range("G8").Value = "We will need to research how to sufficiently create a server that will support a database"
$ws.Range("H8").Value = "Medium "
$ws.Range("I8").Value = "Low "

# ---------------------------------------------------------------------------
# Apply formatting to the new rows by copying from similarly-styled existing
# cells, then fix up row heights / column widths / selection.
# ---------------------------------------------------------------------------

# C6, C8 use the "risk" cell style with border + dark font + fill-applied variant
$ws.Range("C4").Copy()
$ws.Range("C6").PasteSpecial(-4122)
$ws.Range("C8").PasteSpecial(-4122)

# C7 uses the plain bordered style (like header cells)
$ws.Range("C2").Copy()
$ws.Range("C7").PasteSpecial(-4122)

# D6, H7, H8 use the "Medium" style (like E3)
$ws.Range("E3").Copy()
$ws.Range("D6").PasteSpecial(-4122)
$ws.Range("H7").PasteSpecial(-4122)
$ws.Range("H8").PasteSpecial(-4122)

# E6, E7, E8 use the "Medium" style too (like H3)
$ws.Range("H3").Copy()
$ws.Range("E6").PasteSpecial(-4122)
$ws.Range("E7").PasteSpecial(-4122)
$ws.Range("E8").PasteSpecial(-4122)

# D7, D8 use the "High" style (like D5)
$ws.Range("D5").Copy()
$ws.Range("D7").PasteSpecial(-4122)
$ws.Range("D8").PasteSpecial(-4122)

# F6:G8 use the plain bordered style (no fill) like header row cells
$ws.Range("C2").Copy()
$ws.Range("F6:G8").PasteSpecial(-4122)

# H6, I6, I7, I8 use the green "after mitigation" style (like H4/I4)
$ws.Range("H4").Copy()
$ws.Range("H6").PasteSpecial(-4122)
$ws.Range("I6").PasteSpecial(-4122)
$ws.Range("I7").PasteSpecial(-4122)
$ws.Range("I8").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Row heights for the new rows
$ws.Rows.Item(6).RowHeight = 60
$ws.Rows.Item(7).RowHeight = 75
$ws.Rows.Item(8).RowHeight = 60

# Column widths
$ws.Columns.Item(3).ColumnWidth = 15.15
$ws.Columns.Item(7).ColumnWidth = 21.15

# Selection / active cell to match the saved view state
[void]$ws.Range("C2:I8").Select()
